# Update cryptos list (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep Price/Volume columns as plain text so values like "21.00" or
# "0.00001134" are not re-interpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.284.22"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.874.81"
$ws.Range("E3").Value = "  +4.65%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "312.52"
$ws.Range("E5").Value = "  +2.33%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.14%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5063"
$ws.Range("E7").Value = "  +2.26%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3943"
$ws.Range("E8").Value = "  +2.35%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.09586"
$ws.Range("E9").Value = "  +4.00%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "1.145"
$ws.Range("E10").Value = "  +5.21%  "

# Row 11 - OKB
$ws.Range("D11").Value = "40.88"
$ws.Range("E11").Value = "  +0.99%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "6.493"
$ws.Range("E12").Value = "  +3.86%  "

# Row 13 - Solana
$ws.Range("D13").Value = "21.00"
$ws.Range("E13").Value = "  +3.05%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.878.17"
$ws.Range("E14").Value = "  +4.95%  "

# Row 15 - now BinanceUSD (was Chainlink)
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16 - now Chainlink (was BinanceUSD)
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.445"
$ws.Range("E16").Value = "  +4.57%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").Value = "  +3.01%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "93.08"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.06597"
$ws.Range("E19").Value = "  +1.11%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "17.62"
$ws.Range("E20").Value = "  +3.85%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.01%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.191"
$ws.Range("E22").Value = "  +5.17%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "28.330.88"
$ws.Range("E23").Value = "  +2.45%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +3.80%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  +3.72%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "2.571"
$ws.Range("E26").Value = "  +8.14%  "

# Row 27 - WrappedliquidstakedEther2.0
$ws.Range("D27").Value = "2.095.35"
$ws.Range("E27").Value = "  +4.96%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "21.25"
$ws.Range("E28").Value = "  +4.71%  "

# Row 29 - Monero
$ws.Range("D29").Value = "159.03"
$ws.Range("E29").Value = "  +1.34%  "

# Row 30 - BitcoinCash
$ws.Range("E30").Value = "  +1.32%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  +0.66%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +2.23%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.642"
$ws.Range("E33").Value = "  +2.91%  "

# Row 34 - HuobiToken
$ws.Range("D34").Value = "3.629"
$ws.Range("E34").Value = "  +0.68%  "

# Row 35 - FraxShare
$ws.Range("D35").Value = "9.558"
$ws.Range("E35").Value = "  +8.48%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "0.06713"
$ws.Range("E36").Value = "  -0.84%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.02390"
$ws.Range("E37").Value = "  +4.43%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "0.2196"
$ws.Range("E38").Value = "  +3.87%  "

# Row 39 - Aptos
$ws.Range("E39").Value = "  +1.88%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "0.6366"
$ws.Range("E40").Value = "  +4.42%  "

# Row 41 - InternetComputer(DFINITY)
$ws.Range("D41").Value = "5.001"
$ws.Range("E41").Value = "  +2.26%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "1.186"
$ws.Range("E42").Value = "  +4.14%  "

# Row 43 - Frax
$ws.Range("E43").Value = "  +0.06%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "13.54"
$ws.Range("E44").Value = "  +5.00%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "0.5996"
$ws.Range("E45").Value = "  +3.14%  "

# Row 46 - PancakeSwap
$ws.Range("D46").Value = "3.661"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47 - WEMIXTOKEN
$ws.Range("D47").Value = "1.271"
$ws.Range("E47").Value = "  +0.44%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "2.005"
$ws.Range("E48").Value = "  +4.70%  "

# Row 49 - Quant
$ws.Range("D49").Value = "124.26"
$ws.Range("E49").Value = "  +1.32%  "

# Row 50 - EOS
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").Value = "  +2.86%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +2.65%  "
